$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows (14/05/2021 .. 27/05/2021), continuing the existing table
# that ends at row 255 (13/05/2021). Values taken from the source update.
$newData = @(
    @(44330, 0, 2, 67.43088334457181),
    @(44331, 0, 2, 67.43088334457181),
    @(44332, 0, 1, 33.71544167228591),
    @(44333, 1, 2, 67.43088334457181),
    @(44334, 0, 1, 33.71544167228591),
    @(44335, 0, 1, 33.71544167228591),
    @(44336, 0, 1, 33.71544167228591),
    @(44337, 0, 1, 33.71544167228591),
    @(44338, 0, 1, 33.71544167228591),
    @(44339, 0, 1, 33.71544167228591),
    @(44340, 0, 0, 0),
    @(44341, 0, 0, 0),
    @(44342, 0, 0, 0),
    @(44343, 0, 0, 0)
)

$lastRow = 255
$startRow = $lastRow + 1

# Copy the formatting of the last existing row down onto the new rows so the
# date column keeps its style (s="2", the centred date number format), then
# overwrite the values.
$srcRow = $ws.Range("A" + $lastRow + ":D" + $lastRow)
$endRow = $startRow + $newData.Count - 1
$destRange = $ws.Range("A" + $startRow + ":D" + $endRow)
$srcRow.Copy()
$destRange.PasteSpecial(-4122)

for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $startRow + $i
    $row = $newData[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
